$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values in row 1 for columns P and Q, matching the
# existing bold/centered/bordered header style already used by O1.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# For rows 2 through 25, swap values in columns I/K and M/O, and add new
# columns P and Q with value 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I
    $ws.Cells.Item($r, 11).Value = 1  # K
    $ws.Cells.Item($r, 13).Value = 2  # M
    $ws.Cells.Item($r, 15).Value = 1  # O
    $ws.Cells.Item($r, 16).Value = 2  # P
    $ws.Cells.Item($r, 17).Value = 2  # Q
}
